$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.491.48'
$ws.Range('E2').Value = '  -1.03%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.913.46'
$ws.Range('E3').Value = '  -1.48%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.98'
$ws.Range('E5').Value = '  -1.25%  '

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4783'
$ws.Range('E7').Value = '  -2.22%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2843'
$ws.Range('E8').Value = '  -3.77%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06698'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.82'
$ws.Range('E10').Value = '  -3.18%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '101.93'
$ws.Range('E11').Value = '  -4.43%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.920.69'
$ws.Range('E12').Value = '  -1.38%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07694'
$ws.Range('E13').Value = '  -0.27%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.212'
$ws.Range('E14').Value = '  -2.52%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6720'
$ws.Range('E15').Value = '  -3.81%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '267.55'
$ws.Range('E16').Value = '  -3.44%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.492.09'
$ws.Range('E17').Value = '  -1.04%  '

$ws.Range('E18').Value = '  +0.04%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007461'
$ws.Range('E19').Value = '  -3.37%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.67'
$ws.Range('E20').Value = '  -3.27%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.399'
$ws.Range('E21').Value = '  -1.18%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.298'
$ws.Range('E23').Value = '  -3.27%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.369'
$ws.Range('E24').Value = '  -3.61%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '167.22'
$ws.Range('E25').Value = '  -0.27%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '19.17'
$ws.Range('E26').Value = '  -2.44%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.063'
$ws.Range('E27').Value = '  -4.18%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.383'
$ws.Range('E28').Value = '  -0.26%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.669'
$ws.Range('E29').Value = '  +2.54%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09993'
$ws.Range('E30').Value = '  -4.47%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.510'
$ws.Range('E31').Value = '  -2.55%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.239'
$ws.Range('E32').Value = '  -2.81%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04728'
$ws.Range('E33').Value = '  -2.45%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7271'
$ws.Range('E34').Value = '  -3.25%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.108'
$ws.Range('E35').Value = '  -4.28%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.713'
$ws.Range('E36').Value = '  -0.46%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01914'
$ws.Range('E37').Value = '  -3.91%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.619'
$ws.Range('E38').Value = '  -1.53%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.227'
$ws.Range('E39').Value = '  -3.84%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '74.86'
$ws.Range('E40').Value = '  -4.23%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.970'
$ws.Range('E41').Value = '  -5.95%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8625'
$ws.Range('E42').Value = '  -4.81%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '105.69'
$ws.Range('E43').Value = '  -2.17%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4258'
$ws.Range('E44').Value = '  -3.16%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9997'
$ws.Range('E45').Value = '  +0.12%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.406'
$ws.Range('E46').Value = '  -4.49%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1203'
$ws.Range('E47').Value = '  -3.23%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '914.36'
$ws.Range('E48').Value = '  -7.29%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '34.71'
$ws.Range('E49').Value = '  -3.44%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.747'
$ws.Range('E50').Value = '  -5.81%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05746'
$ws.Range('E51').Value = '  +0.31%  '
